$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10: set the URL cell first so it claims the earlier shared-string
# index, then the label cell, matching the order new strings were appended
# in the target workbook.
$ws.Range("C10").Value = "https://github.com/nguyentienminh07102004/product-management/commit/9516b17e3bbe0f9823dc17f0bdf0c65b8008935b"
$ws.Range("B10").Value = "3.2. Hiển thị thông báo sau khi xóa"

# Turn the URL cell into a real hyperlink (adds the relationship + applies
# the Hyperlink cell style).
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/nguyentienminh07102004/product-management/commit/9516b17e3bbe0f9823dc17f0bdf0c65b8008935b") | Out-Null
$ws.Range("C10").Style = "Hyperlink"

# Move selection down to C11, as in the saved workbook.
$ws.Range("C11").Select() | Out-Null
